$wb = $excel.ActiveWorkbook

# --- Sheet ALC: update cached leve-profit values ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4417.647
$ws.Range("I64").Value = 3724.25
$ws.Range("K64").Value = 3724.25
$ws.Range("M64").Value = -3476.25
$ws.Range("H67").Value = 4417.647
$ws.Range("I67").Value = 3724.25
$ws.Range("K67").Value = 3724.25
$ws.Range("M67").Value = -2866.25
$ws.Range("H82").Value = 10577.333
$ws.Range("I82").Value = 1503.5
$ws.Range("J82").Value = 28725
$ws.Range("K82").Value = 4510.5
$ws.Range("L82").Value = 86175
$ws.Range("M82").Value = -4104.5
$ws.Range("N82").Value = -86987
$ws.Range("H85").Value = 10577.333
$ws.Range("I85").Value = 1503.5
$ws.Range("J85").Value = 28725
$ws.Range("K85").Value = 4510.5
$ws.Range("L85").Value = 86175
$ws.Range("M85").Value = -3106.5
$ws.Range("N85").Value = -88983
$ws.Range("H106").Value = 1009.3333
$ws.Range("I106").Value = 885.5
$ws.Range("J106").Value = 2000
$ws.Range("K106").Value = 885.5
$ws.Range("L106").Value = 2000
$ws.Range("M106").Value = -254.5
$ws.Range("N106").Value = -3262
$ws.Range("H112").Value = 2035.4517
$ws.Range("J112").Value = 2157.0715
$ws.Range("L112").Value = 6471.2145
$ws.Range("N112").Value = -8687.2145
$ws.Range("H125").Value = 1022
$ws.Range("I125").Value = 779.3333
$ws.Range("J125").Value = 1750
$ws.Range("K125").Value = 7013.9997
$ws.Range("L125").Value = 15750
$ws.Range("M125").Value = -4553.9997
$ws.Range("N125").Value = -20670
$ws.Range("H129").Value = 976.73334
$ws.Range("J129").Value = 1121.0857
$ws.Range("L129").Value = 3363.2571
$ws.Range("N129").Value = -13363.2571

# --- Sheet ARM: update cached leve-profit values ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1199.9166
$ws.Range("I2").Value = 1016.5
$ws.Range("J2").Value = 1383.3334
$ws.Range("K2").Value = 1016.5
$ws.Range("L2").Value = 1383.3334
$ws.Range("M2").Value = -903.5
$ws.Range("N2").Value = -1609.3334
$ws.Range("H63").Value = 8381.8125
$ws.Range("I63").Value = 7301.1113
$ws.Range("J63").Value = 9771.286
$ws.Range("K63").Value = 7301.1113
$ws.Range("L63").Value = 9771.286
$ws.Range("M63").Value = -6615.1113
$ws.Range("N63").Value = -11143.286
$ws.Range("H66").Value = 8381.8125
$ws.Range("I66").Value = 7301.1113
$ws.Range("J66").Value = 9771.286
$ws.Range("K66").Value = 36505.5565
$ws.Range("L66").Value = 48856.43
$ws.Range("M66").Value = -33073.5565
$ws.Range("N66").Value = -55720.43
$ws.Range("H116").Value = 1199.9166
$ws.Range("I116").Value = 1016.5
$ws.Range("J116").Value = 1383.3334
$ws.Range("K116").Value = 1016.5
$ws.Range("L116").Value = 1383.3334
$ws.Range("M116").Value = 1277.5
$ws.Range("N116").Value = -5971.3334
$ws.Range("H122").Value = 2041.5217
$ws.Range("I122").Value = 2201.6875
$ws.Range("J122").Value = 1675.4286
$ws.Range("K122").Value = 6605.0625
$ws.Range("L122").Value = 5026.2858
$ws.Range("M122").Value = -4155.0625
$ws.Range("N122").Value = -9926.2858

# --- Sheet BSM: update cached leve-profit values ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1199.9166
$ws.Range("I3").Value = 1016.5
$ws.Range("J3").Value = 1383.3334
$ws.Range("K3").Value = 1016.5
$ws.Range("L3").Value = 1383.3334
$ws.Range("M3").Value = -902.5
$ws.Range("N3").Value = -1611.3334
$ws.Range("H74").Value = 65390
$ws.Range("J74").Value = 65390
$ws.Range("L74").Value = 65390
$ws.Range("N74").Value = -67262
$ws.Range("H77").Value = 65390
$ws.Range("J77").Value = 65390
$ws.Range("L77").Value = 196170
$ws.Range("N77").Value = -205530
$ws.Range("H105").Value = 3337.7646
$ws.Range("I105").Value = 2924.4285
$ws.Range("K105").Value = 2924.4285
$ws.Range("M105").Value = -1177.4285
$ws.Range("H107").Value = 29582.475
$ws.Range("I107").Value = 44843.082
$ws.Range("J107").Value = 3421.4285
$ws.Range("K107").Value = 44843.082
$ws.Range("L107").Value = 3421.4285
$ws.Range("M107").Value = -42923.082
$ws.Range("N107").Value = -7261.4285

# --- Sheet CRP: update cached leve-profit values ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2845.652
$ws.Range("I31").Value = 2053.697
$ws.Range("J31").Value = 4856
$ws.Range("K31").Value = 2053.697
$ws.Range("L31").Value = 4856
$ws.Range("M31").Value = -1758.697
$ws.Range("N31").Value = -5446
$ws.Range("H34").Value = 2845.652
$ws.Range("I34").Value = 2053.697
$ws.Range("J34").Value = 4856
$ws.Range("K34").Value = 2053.697
$ws.Range("L34").Value = 4856
$ws.Range("M34").Value = -1851.697
$ws.Range("N34").Value = -5260

# --- Sheet CUL: update cached leve-profit values ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1669.5
$ws.Range("I109").Value = 661.9286
$ws.Range("K109").Value = 1985.7858
$ws.Range("M109").Value = -945.7857999999999
$ws.Range("H124").Value = 7859.8
$ws.Range("I124").Value = 1433
$ws.Range("J124").Value = 17500
$ws.Range("K124").Value = 4299
$ws.Range("L124").Value = 52500
$ws.Range("M124").Value = 611
$ws.Range("N124").Value = -62320

# --- Sheet GSM: update cached leve-profit values ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6350.4736
$ws.Range("I70").Value = 5350
$ws.Range("J70").Value = 6617.2666
$ws.Range("K70").Value = 5350
$ws.Range("L70").Value = 6617.2666
$ws.Range("M70").Value = -5080
$ws.Range("N70").Value = -7157.2666
$ws.Range("H73").Value = 6350.4736
$ws.Range("I73").Value = 5350
$ws.Range("J73").Value = 6617.2666
$ws.Range("K73").Value = 5350
$ws.Range("L73").Value = 6617.2666
$ws.Range("M73").Value = -4414
$ws.Range("N73").Value = -8489.266599999999

# --- Sheet LTW: update cached leve-profit values ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 15181.267
$ws.Range("J61").Value = 2924.5
$ws.Range("L61").Value = 2924.5
$ws.Range("N61").Value = -3328.5
$ws.Range("H113").Value = 15181.267
$ws.Range("J113").Value = 2924.5
$ws.Range("L113").Value = 2924.5
$ws.Range("N113").Value = -7264.5

# --- Sheet WVR: update cached leve-profit values ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 202960
$ws.Range("I81").Value = 334933.34
$ws.Range("K81").Value = 669866.6800000001
$ws.Range("M81").Value = -668805.6800000001
$ws.Range("H84").Value = 202960
$ws.Range("I84").Value = 334933.34
$ws.Range("K84").Value = 3349333.4
$ws.Range("M84").Value = -3344029.4
$ws.Range("H133").Value = 44482
$ws.Range("J133").Value = 44482
$ws.Range("L133").Value = 44482
$ws.Range("N133").Value = -54602
